$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date column (C2:C6) from 2023-10-09 (45208) to 2023-10-13 (45212)
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
